$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 46071
$ws.Cells.Item(2, 2).Value = 10.333
$ws.Cells.Item(2, 3).Value = 0.051
$ws.Cells.Item(3, 1).Value = 46071.01041666666
$ws.Cells.Item(3, 2).Value = 5.201
$ws.Cells.Item(3, 3).Value = 0.076
$ws.Cells.Item(4, 1).Value = 46071.02083333334
$ws.Cells.Item(4, 2).Value = 22.104
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(5, 1).Value = 46071.03125
$ws.Cells.Item(5, 2).Value = 28.747
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(6, 1).Value = 46071.04166666666
$ws.Cells.Item(6, 2).Value = 7.88
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(7, 1).Value = 46071.05208333334
$ws.Cells.Item(7, 2).Value = 9.317
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(8, 1).Value = 46071.0625
$ws.Cells.Item(8, 2).Value = 10.206
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(9, 1).Value = 46071.07291666666
$ws.Cells.Item(9, 2).Value = 22.984
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(10, 1).Value = 46071.08333333334
$ws.Cells.Item(10, 2).Value = 16.757
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(11, 1).Value = 46071.09375
$ws.Cells.Item(11, 2).Value = 7.853
$ws.Cells.Item(11, 3).Value = 0.342
$ws.Cells.Item(12, 1).Value = 46071.10416666666
$ws.Cells.Item(12, 2).Value = 29.304
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(13, 1).Value = 46071.11458333334
$ws.Cells.Item(13, 2).Value = 46.916
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(14, 1).Value = 46071.125
$ws.Cells.Item(14, 2).Value = 17.701
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(15, 1).Value = 46071.13541666666
$ws.Cells.Item(15, 2).Value = 10.935
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(16, 1).Value = 46071.14583333334
$ws.Cells.Item(16, 2).Value = 26.933
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(17, 1).Value = 46071.15625
$ws.Cells.Item(17, 2).Value = 41.368
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(18, 1).Value = 46071.16666666666
$ws.Cells.Item(18, 2).Value = 2.442
$ws.Cells.Item(18, 3).Value = 1.145
$ws.Cells.Item(19, 1).Value = 46071.17708333334
$ws.Cells.Item(19, 2).Value = 9.074999999999999
$ws.Cells.Item(19, 3).Value = 0.288
$ws.Cells.Item(20, 1).Value = 46071.1875
$ws.Cells.Item(20, 2).Value = 24.412
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(21, 1).Value = 46071.19791666666
$ws.Cells.Item(21, 2).Value = 7.959
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(22, 1).Value = 46071.20833333334
$ws.Cells.Item(22, 2).Value = 0.005
$ws.Cells.Item(22, 3).Value = 18.52
$ws.Cells.Item(23, 1).Value = 46071.21875
$ws.Cells.Item(23, 2).Value = 0.012
$ws.Cells.Item(23, 3).Value = 9.137
$ws.Cells.Item(24, 1).Value = 46071.22916666666
$ws.Cells.Item(24, 2).Value = 0.919
$ws.Cells.Item(24, 3).Value = 4.172
$ws.Cells.Item(25, 1).Value = 46071.23958333334
$ws.Cells.Item(25, 2).Value = 9.045999999999999
$ws.Cells.Item(25, 3).Value = 0.372
$ws.Cells.Item(26, 1).Value = 46071.25
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(26, 3).Value = 10.637
$ws.Cells.Item(27, 1).Value = 46071.26041666666
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(27, 3).Value = 17.986
$ws.Cells.Item(28, 1).Value = 46071.27083333334
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(28, 3).Value = 17.288
$ws.Cells.Item(29, 1).Value = 46071.28125
$ws.Cells.Item(29, 2).Value = 1.287
$ws.Cells.Item(29, 3).Value = 2.822
$ws.Cells.Item(30, 1).Value = 46071.29166666666
$ws.Cells.Item(30, 2).Value = 20.425
$ws.Cells.Item(30, 3).Value = 0.423
$ws.Cells.Item(31, 1).Value = 46071.3125
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 0
